$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Compounds_of_interest")

$ws.Range("A2").Value = "my_compound_1"
$ws.Range("A3").Value = "my_compound_2"
$ws.Range("A4").Value = "my_compound_3"
$ws.Range("A5").Value = "my_compound_4"
$ws.Range("A6").Value = "my_compound_5"

$ws.Activate()
$ws.Range("A2").Select()
